$d = $word.ActiveDocument

# --- Change 1: collapse the "deltaTime" paragraph's split runs (with the
# spell-check proofErr markers around it) into a single run of the same
# text. A Find/Replace across the whole paragraph text rewrites the
# matched range as one run and drops the now-irrelevant proofErr markers.
$old1 = "Created the log. Continued working on the pour prototype. Spent a few hours attempting multiple solutions to create the cursor system I ended up with. It has click and hold functionality right now, need to work on adding a click function (using deltaTime to space it out, probably). The keg, tray, and cup are all interactable."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $old1, 2) | Out-Null

# --- Change 2: split "Differentiated click and hold !" into two runs:
# "Differentiated click and hold" followed by a new run containing
# ". Cup now moves to tray when clicked." ---
#
# A plain Range.InsertAfter right after a Find/Replace would type the new
# text into the *same* run (adjacent runs with identical formatting get
# coalesced back together), so instead the paragraph is temporarily split
# in two, each half gets its own text, and then the leading half's
# paragraph mark is deleted to rejoin them. Deleting a paragraph mark
# merges its content forward into the *next* paragraph while that next
# paragraph's identity (w14:paraId/rsids) is what survives - so the new
# "prefix" paragraph is inserted *before* the original one (whose text is
# first replaced with just the new suffix), keeping the original
# paragraph's identity on the final, merged paragraph while still
# producing two independent <w:r> runs in the right order.

# 1) Replace the old sentence in place (same run/paragraph) with just the
#    new suffix text.
$d.Content.Find.Execute("Differentiated click and hold !", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ". Cup now moves to tray when clicked.", 2) | Out-Null

# 2) Insert a fresh, blank paragraph immediately before it.
$r = $d.Content
$r.Find.Execute(". Cup now moves to tray when clicked.", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0) | Out-Null
$r.Collapse(1)
$r.InsertParagraphBefore()

# 3) Fill the new (now-empty) paragraph with the "prefix" text.
$prefixPara = $d.Paragraphs(4)
$suffixPara = $d.Paragraphs(5)
$prefixPara.Range.InsertBefore("Differentiated click and hold")

# 4) Delete the prefix paragraph's own paragraph mark so it merges forward
#    into the (identity-preserving) suffix paragraph, leaving two runs.
$markStart = $prefixPara.Range.End - 1
$markRange = $d.Range($markStart, $markStart + 1)
$markRange.Delete()
